$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (DR), shifting DR..Guidelines one
# column to the right, and give it the new required-field header.
$ws.Range("H1").EntireColumn.Insert()
$ws.Cells.Item(1, 8).Value = "Insertion currency"

# Match the column width the new "Insertion currency" column ended up with.
$ws.Columns.Item(8).ColumnWidth = 16.7

# Move the active selection, as recorded after the edit.
$ws.Range("H4").Select()
